$wb = $excel.ActiveWorkbook

# Update last-updated timestamp on Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 03:18 PM"

# Update Stock List sheet data
$ws = $wb.Worksheets.Item("Stock List")
$ws.Range("B2").Value = "MIDWESTLTD"
$ws.Range("C2").Value = "MIDWESTLTD"
$ws.Range("D2").Value = 1117.2
$ws.Range("E2").Value = -1.4032
$ws.Range("H2").Value = 4039.8864
$ws.Range("B3").Value = "CAPTRU-RE1"
$ws.Range("C3").Value = "CAPTRU-RE1"
$ws.Range("D3").Value = 5.67
$ws.Range("E3").Value = -11.9565
$ws.Range("H3").Value = 0
$ws.Range("B4").Value = "NIFTYCASE"
$ws.Range("C4").Value = "NIFTYCASE"
$ws.Range("D4").Value = 10.19
$ws.Range("E4").Value = -0.5854
$ws.Range("H4").Value = 0
$ws.Range("B5").Value = "MOMENTUM30"
$ws.Range("C5").Value = "MOMENTUM30"
$ws.Range("D5").Value = 31.54
$ws.Range("E5").Value = -0.6614
$ws.Range("H5").Value = 0
$ws.Range("B6").Value = "CANHLIFE"
$ws.Range("C6").Value = "CANHLIFE"
$ws.Range("D6").Value = 118.46
$ws.Range("E6").Value = 0.6286
$ws.Range("H6").Value = 11253.7
$ws.Range("B7").Value = "FLEXIADD"
$ws.Range("C7").Value = "FLEXIADD"
$ws.Range("D7").Value = 10.64
$ws.Range("E7").Value = -1.0233
$ws.Range("H7").Value = 0
$ws.Range("B8").Value = "MOENERGY"
$ws.Range("C8").Value = "MOENERGY"
$ws.Range("D8").Value = 36.3
$ws.Range("E8").Value = -0.6568000000000001
$ws.Range("H8").Value = 0
$ws.Range("B9").Value = "MONIFTY100"
$ws.Range("C9").Value = "MONIFTY100"
$ws.Range("D9").Value = 26.49
$ws.Range("E9").Value = 0.3409
$ws.Range("H9").Value = 0
$ws.Range("B10").Value = "RUBICON"
$ws.Range("C10").Value = "RUBICON"
$ws.Range("D10").Value = 652.65
$ws.Range("E10").Value = -0.1453
$ws.Range("H10").Value = 10752.4289
$ws.Range("B11").Value = "CRAMC"
$ws.Range("C11").Value = "CRAMC"
$ws.Range("D11").Value = 317.2
$ws.Range("E11").Value = 2.3226
$ws.Range("H11").Value = 6325.5208
$ws.Range("B12").Value = "LGEINDIA"
$ws.Range("C12").Value = "LGEINDIA"
$ws.Range("D12").Value = 1633.4
$ws.Range("E12").Value = -0.946
$ws.Range("H12").Value = 110870.6825
$ws.Range("B13").Value = "TATACAP"
$ws.Range("C13").Value = "TATACAP"
$ws.Range("D13").Value = 329.3
$ws.Range("E13").Value = 0.1521
$ws.Range("H13").Value = 139783.5374
$ws.Range("B14").Value = "ELIQUID"
$ws.Range("C14").Value = "ELIQUID"
$ws.Range("D14").Value = 1004.85
$ws.Range("E14").Value = 0.0408
$ws.Range("H14").Value = 0
$ws.Range("B15").Value = "WEWORK"
$ws.Range("C15").Value = "WEWORK"
$ws.Range("D15").Value = 632.15
$ws.Range("E15").Value = -2.4008
$ws.Range("H15").Value = 8472.2803
$ws.Range("B16").Value = "GROWWRLTY"
$ws.Range("C16").Value = "GROWWRLTY"
$ws.Range("D16").Value = 10.8
$ws.Range("E16").Value = -0.4608
$ws.Range("H16").Value = 0
$ws.Range("B17").Value = "ADVANCE"
$ws.Range("C17").Value = "ADVANCE"
$ws.Range("D17").Value = 130.05
$ws.Range("E17").Value = -5.2666
$ws.Range("H17").Value = 836.0358
$ws.Range("B18").Value = "OMFREIGHT"
$ws.Range("C18").Value = "OMFREIGHT"
$ws.Range("D18").Value = 88.90000000000001
$ws.Range("E18").Value = -0.5926
$ws.Range("H18").Value = 299.3747
$ws.Range("B19").Value = "GLOTTIS"
$ws.Range("C19").Value = "GLOTTIS"
$ws.Range("D19").Value = 72.73999999999999
$ws.Range("E19").Value = -0.8587
$ws.Range("H19").Value = 672.1394
$ws.Range("B20").Value = "FABTECH"
$ws.Range("C20").Value = "FABTECH"
$ws.Range("D20").Value = 237.72
$ws.Range("E20").Value = 0.4734
$ws.Range("H20").Value = 1056.6843
$ws.Range("B21").Value = "PACEDIGITK"
$ws.Range("C21").Value = "PACEDIGITK"
$ws.Range("D21").Value = 218.85
$ws.Range("E21").Value = 0.1327
$ws.Range("H21").Value = 4723.9063
$ws.Range("B22").Value = "JAINREC"
$ws.Range("C22").Value = "JAINREC"
$ws.Range("D22").Value = 377.25
$ws.Range("E22").Value = 1.2208
$ws.Range("H22").Value = 13018.3623
$ws.Range("B23").Value = "EPACKPEB"
$ws.Range("C23").Value = "EPACKPEB"
$ws.Range("D23").Value = 301.45
$ws.Range("E23").Value = 1.979
$ws.Range("H23").Value = 3028.1254
$ws.Range("B24").Value = "BMWVENTLTD"
$ws.Range("C24").Value = "BMWVENTLTD"
$ws.Range("D24").Value = 69.25
$ws.Range("E24").Value = 0
$ws.Range("H24").Value = 600.5014
$ws.Range("B25").Value = "STYL"
$ws.Range("C25").Value = "STYL"
$ws.Range("D25").Value = 372.4
$ws.Range("E25").Value = -0.8388
$ws.Range("H25").Value = 6025.649
$ws.Range("B26").Value = "JARO"
$ws.Range("C26").Value = "JARO"
$ws.Range("D26").Value = 621.5
$ws.Range("E26").Value = -1.4821
$ws.Range("H26").Value = 1377.0134
$ws.Range("B27").Value = "SOLARWORLD"
$ws.Range("C27").Value = "SOLARWORLD"
$ws.Range("D27").Value = 309.1
$ws.Range("E27").Value = -0.6269
$ws.Range("H27").Value = 2679.0517
$ws.Range("B28").Value = "ARSSBL"
$ws.Range("C28").Value = "ARSSBL"
$ws.Range("D28").Value = 537.3
$ws.Range("E28").Value = 4.7266
$ws.Range("H28").Value = 3370.2277
$ws.Range("B29").Value = "GANESHCP"
$ws.Range("C29").Value = "GANESHCP"
$ws.Range("D29").Value = 274.4
$ws.Range("E29").Value = -2.7984
$ws.Range("H29").Value = 1108.9312
$ws.Range("B30").Value = "ATLANTAELE"
$ws.Range("C30").Value = "ATLANTAELE"
$ws.Range("D30").Value = 1003.05
$ws.Range("E30").Value = -1.7436
$ws.Range("H30").Value = 7713.116
$ws.Range("B31").Value = "GKENERGY"
$ws.Range("C31").Value = "GKENERGY"
$ws.Range("D31").Value = 213.85
$ws.Range("E31").Value = -0.7933
$ws.Range("H31").Value = 4337.2472
$ws.Range("B32").Value = "SAATVIKGL"
$ws.Range("C32").Value = "SAATVIKGL"
$ws.Range("D32").Value = 528.2
$ws.Range("E32").Value = -1.3079
$ws.Range("H32").Value = 6713.6863
$ws.Range("B33").Value = "IVALUE"
$ws.Range("C33").Value = "IVALUE"
$ws.Range("D33").Value = 281.45
$ws.Range("E33").Value = -0.3364
$ws.Range("H33").Value = 1506.8799
$ws.Range("B34").Value = "VMSTMT"
$ws.Range("C34").Value = "VMSTMT"
$ws.Range("D34").Value = 70.03
$ws.Range("E34").Value = -0.9056
$ws.Range("H34").Value = 347.5674
$ws.Range("B35").Value = "EUROPRATIK"
$ws.Range("C35").Value = "EUROPRATIK"
$ws.Range("D35").Value = 321.75
$ws.Range("E35").Value = 0.8147
$ws.Range("H35").Value = 3288.285
$ws.Range("B36").Value = "SHRINGARMS"
$ws.Range("C36").Value = "SHRINGARMS"
$ws.Range("D36").Value = 229.31
$ws.Range("E36").Value = -1.2616
$ws.Range("H36").Value = 2211.284
$ws.Range("B37").Value = "DEVX"
$ws.Range("C37").Value = "DEVX"
$ws.Range("D37").Value = 44.53
$ws.Range("E37").Value = -0.3803
$ws.Range("H37").Value = 401.605
$ws.Range("B38").Value = "URBANCO"
$ws.Range("C38").Value = "URBANCO"
$ws.Range("D38").Value = 148.9
$ws.Range("E38").Value = -2.0459
$ws.Range("H38").Value = 21380.5798
$ws.Range("B39").Value = "SML100CASE"
$ws.Range("C39").Value = "SML100CASE"
$ws.Range("D39").Value = 10.36
$ws.Range("E39").Value = -0.7663
$ws.Range("H39").Value = 0
$ws.Range("B40").Value = "AONEGOLD"
$ws.Range("C40").Value = "AONEGOLD"
$ws.Range("D40").Value = 11.28
$ws.Range("E40").Value = -0.2653
$ws.Range("H40").Value = 0
$ws.Range("B41").Value = "ELM250"
$ws.Range("C41").Value = "ELM250"
$ws.Range("D41").Value = 16.72
$ws.Range("E41").Value = 0.1797
$ws.Range("H41").Value = 0
$ws.Range("B42").Value = "AMANTA"
$ws.Range("C42").Value = "AMANTA"
$ws.Range("D42").Value = 122.52
$ws.Range("E42").Value = 1.407
$ws.Range("H42").Value = 475.7372
$ws.Range("B43").Value = "CPEDU"
$ws.Range("C43").Value = "CPEDU"
$ws.Range("D43").Value = 315.9
$ws.Range("E43").Value = 1.8539
$ws.Range("H43").Value = 574.7148999999999
$ws.Range("B44").Value = "AHCL"
$ws.Range("C44").Value = "AHCL"
$ws.Range("D44").Value = 139.27
$ws.Range("E44").Value = 3.1706
$ws.Range("H44").Value = 740.2409
$ws.Range("B45").Value = "STLNETWORK"
$ws.Range("C45").Value = "STLNETWORK"
$ws.Range("D45").Value = 26.59
$ws.Range("E45").Value = -0.412
$ws.Range("H45").Value = 1297.3822
$ws.Range("B46").Value = "VIKRAN"
$ws.Range("C46").Value = "VIKRAN"
$ws.Range("D46").Value = 98.05
$ws.Range("E46").Value = -1.783
$ws.Range("H46").Value = 2528.8166
$ws.Range("B47").Value = "MANUFGBEES"
$ws.Range("C47").Value = "MANUFGBEES"
$ws.Range("D47").Value = 151.77
$ws.Range("E47").Value = -1.011
$ws.Range("H47").Value = 0
$ws.Range("B48").Value = "MEIL"
$ws.Range("C48").Value = "MEIL"
$ws.Range("D48").Value = 461.15
$ws.Range("E48").Value = -0.7319
$ws.Range("H48").Value = 1274.1632
$ws.Range("B49").Value = "GROWWNXT50"
$ws.Range("C49").Value = "GROWWNXT50"
$ws.Range("D49").Value = 70.29000000000001
$ws.Range("E49").Value = -0.4109
$ws.Range("H49").Value = 0
$ws.Range("B50").Value = "SHREEJISPG"
$ws.Range("C50").Value = "SHREEJISPG"
$ws.Range("D50").Value = 270.05
$ws.Range("E50").Value = -0.7899
$ws.Range("H50").Value = 4399.6074
$ws.Range("B51").Value = "GEMAROMA"
$ws.Range("C51").Value = "GEMAROMA"
$ws.Range("D51").Value = 219.52
$ws.Range("E51").Value = -0.876
$ws.Range("H51").Value = 1146.7097
$ws.Range("B52").Value = "PATELRMART"
$ws.Range("C52").Value = "PATELRMART"
$ws.Range("D52").Value = 219.31
$ws.Range("E52").Value = -1.0646
$ws.Range("H52").Value = 732.5069999999999
$ws.Range("B53").Value = "VIKRAMSOLR"
$ws.Range("C53").Value = "VIKRAMSOLR"
$ws.Range("D53").Value = 322
$ws.Range("E53").Value = -1.5892
$ws.Range("H53").Value = 11647.2884
$ws.Range("B54").Value = "LTGILTCASE"
$ws.Range("C54").Value = "LTGILTCASE"
$ws.Range("D54").Value = 29.67
$ws.Range("E54").Value = 0.2365
$ws.Range("H54").Value = 0
$ws.Range("B55").Value = "REGAAL"
$ws.Range("C55").Value = "REGAAL"
$ws.Range("D55").Value = 89.13
$ws.Range("E55").Value = -0.8675
$ws.Range("H55").Value = 915.5742
$ws.Range("B56").Value = "BLUESTONE"
$ws.Range("C56").Value = "BLUESTONE"
$ws.Range("D56").Value = 711.95
$ws.Range("E56").Value = 0.1266
$ws.Range("H56").Value = 10773.2539
$ws.Range("B57").Value = "MOSILVER"
$ws.Range("C57").Value = "MOSILVER"
$ws.Range("D57").Value = 145.9
$ws.Range("E57").Value = -1.5054
$ws.Range("H57").Value = 0
$ws.Range("B58").Value = "ALLTIME"
$ws.Range("C58").Value = "ALLTIME"
$ws.Range("D58").Value = 308.75
$ws.Range("E58").Value = 2.66
$ws.Range("H58").Value = 2022.5526
$ws.Range("B59").Value = "JSWCEMENT"
$ws.Range("C59").Value = "JSWCEMENT"
$ws.Range("D59").Value = 134.98
$ws.Range("E59").Value = -0.4793
$ws.Range("H59").Value = 18402.6999
$ws.Range("B60").Value = "SBILIQETF"
$ws.Range("C60").Value = "SBILIQETF"
$ws.Range("D60").Value = 1012.94
$ws.Range("E60").Value = 0.0296
$ws.Range("H60").Value = 0
$ws.Range("B61").Value = "HILINFRA"
$ws.Range("C61").Value = "HILINFRA"
$ws.Range("D61").Value = 77.23
$ws.Range("E61").Value = -0.3998
$ws.Range("H61").Value = 0
$ws.Range("B62").Value = "GROWWPOWER"
$ws.Range("C62").Value = "GROWWPOWER"
$ws.Range("D62").Value = 10.28
$ws.Range("E62").Value = -0.9634
$ws.Range("H62").Value = 0
$ws.Range("B63").Value = "LOTUSDEV"
$ws.Range("C63").Value = "LOTUSDEV"
$ws.Range("D63").Value = 177.82
$ws.Range("E63").Value = 0.3669
$ws.Range("H63").Value = 8690.485000000001
$ws.Range("B64").Value = "MBEL"
$ws.Range("C64").Value = "MBEL"
$ws.Range("D64").Value = 450.2
$ws.Range("E64").Value = -0.7714
$ws.Range("H64").Value = 2572.8126
$ws.Range("B65").Value = "LAXMIINDIA"
$ws.Range("C65").Value = "LAXMIINDIA"
$ws.Range("D65").Value = 145.62
$ws.Range("E65").Value = -1.1942
$ws.Range("H65").Value = 761.1248000000001
$ws.Range("B66").Value = "CPPLUS"
$ws.Range("C66").Value = "CPPLUS"
$ws.Range("D66").Value = 1322.1
$ws.Range("E66").Value = -0.264
$ws.Range("H66").Value = 15497.9053
$ws.Range("B67").Value = "SHANTIGOLD"
$ws.Range("C67").Value = "SHANTIGOLD"
$ws.Range("D67").Value = 241.57
$ws.Range("E67").Value = -1.6409
$ws.Range("H67").Value = 1741.6231
$ws.Range("B68").Value = "MOGOLD"
$ws.Range("C68").Value = "MOGOLD"
$ws.Range("D68").Value = 119.65
$ws.Range("E68").Value = -0.5403
$ws.Range("H68").Value = 0
$ws.Range("B69").Value = "BRIGHOTEL"
$ws.Range("C69").Value = "BRIGHOTEL"
$ws.Range("D69").Value = 82.39
$ws.Range("E69").Value = -0.9855
$ws.Range("H69").Value = 3129.5229
$ws.Range("B70").Value = "INDIQUBE"
$ws.Range("C70").Value = "INDIQUBE"
$ws.Range("D70").Value = 212.64
$ws.Range("E70").Value = -0.7561
$ws.Range("H70").Value = 4465.6847
$ws.Range("B71").Value = "EBGNG"
$ws.Range("C71").Value = "EBGNG"
$ws.Range("D71").Value = 346.65
$ws.Range("E71").Value = 3.2311
$ws.Range("H71").Value = 3952.2092
$ws.Range("B72").Value = "LIQGRWBEES"
$ws.Range("C72").Value = "LIQGRWBEES"
$ws.Range("D72").Value = 1014.74
$ws.Range("E72").Value = 0.0246
$ws.Range("H72").Value = 0
$ws.Range("B73").Value = "CHEMBONDCH"
$ws.Range("C73").Value = "CHEMBONDCH"
$ws.Range("D73").Value = 153.35
$ws.Range("E73").Value = -1.6987
$ws.Range("H73").Value = 412.459
$ws.Range("B74").Value = "GROWWNIFTY"
$ws.Range("C74").Value = "GROWWNIFTY"
$ws.Range("D74").Value = 10.29
$ws.Range("E74").Value = -0.3872
$ws.Range("H74").Value = 0
$ws.Range("B75").Value = "ANTHEM"
$ws.Range("C75").Value = "ANTHEM"
$ws.Range("D75").Value = 702.25
$ws.Range("E75").Value = -0.1209
$ws.Range("H75").Value = 39439.0658
$ws.Range("B76").Value = "QUALITY30"
$ws.Range("C76").Value = "QUALITY30"
$ws.Range("D76").Value = 21.05
$ws.Range("E76").Value = -0.8945
$ws.Range("H76").Value = 0
